$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data for rows 2-33 (columns A-R), reflecting weekly refresh of "Pepino dulce" price records.
$data = @(
  ,@(1, "Agrícola del Norte S.A. de Arica", "Arica y Parinacota", 44748, 15, 100112043, "Pepino dulce", "Cultivar IV Región", "Primera", 250, 17000, 18000, 17500, "`$/bandeja 18 kilos", "Provincia de Limarí", 972, 18, "Hortaliza")
  ,@(1, "Agrícola del Norte S.A. de Arica", "Arica y Parinacota", 44433, 15, 100112043, "Pepino dulce", "Cultivar IV Región", "Segunda", 100, 17000, 18000, 17500, "`$/bandeja 18 kilos", "Provincia de Limarí", 972, 18, "Hortaliza")
  ,@(1, "Agrícola del Norte S.A. de Arica", "Arica y Parinacota", 44433, 15, 100112043, "Pepino dulce", "Cultivar IV Región", "Tercera", 120, 14000, 15000, 14500, "`$/bandeja 18 kilos", "Provincia de Limarí", 806, 18, "Hortaliza")
  ,@(1, "Agrícola del Norte S.A. de Arica", "Arica y Parinacota", 44435, 15, 100112043, "Pepino dulce", "Cultivar IV Región", "Segunda", 100, 17000, 18000, 17500, "`$/bandeja 18 kilos", "Provincia de Limarí", 972, 18, "Hortaliza")
  ,@(1, "Agrícola del Norte S.A. de Arica", "Arica y Parinacota", 44435, 15, 100112043, "Pepino dulce", "Cultivar IV Región", "Tercera", 120, 14000, 15000, 14500, "`$/bandeja 18 kilos", "Provincia de Limarí", 806, 18, "Hortaliza")
  ,@(1, "Agrícola del Norte S.A. de Arica", "Arica y Parinacota", 44405, 15, 100112043, "Pepino dulce", "Cultivar IV Región", "Segunda", 140, 17000, 18000, 17500, "`$/bandeja 18 kilos", "Provincia de Limarí", 972, 18, "Hortaliza")
  ,@(1, "Agrícola del Norte S.A. de Arica", "Arica y Parinacota", 44771, 15, 100112043, "Pepino dulce", "Cultivar XV región", "Primera", 140, 8000, 9000, 8500, "`$/caja 10 kilos", "Región de Arica y Parinacota", 850, 10, "Hortaliza")
  ,@(1, "Agrícola del Norte S.A. de Arica", "Arica y Parinacota", 44742, 15, 100112043, "Pepino dulce", "Cultivar IV Región", "Segunda", 250, 15000, 16000, 15500, "`$/bandeja 18 kilos", "Provincia de Limarí", 861, 18, "Hortaliza")
  ,@(1, "Agrícola del Norte S.A. de Arica", "Arica y Parinacota", 44755, 15, 100112043, "Pepino dulce", "Cultivar IV Región", "Primera", 160, 17000, 18000, 17500, "`$/bandeja 18 kilos", "Provincia de Limarí", 972, 18, "Hortaliza")
  ,@(1, "Agrícola del Norte S.A. de Arica", "Arica y Parinacota", 44526, 15, 100112043, "Pepino dulce", "Cultivar XV región", "Primera", 100, 5000, 5500, 5250, "`$/caja 10 kilos", "Región de Arica y Parinacota", 525, 10, "Hortaliza")
  ,@(1, "Agrícola del Norte S.A. de Arica", "Arica y Parinacota", 44526, 15, 100112043, "Pepino dulce", "Cultivar XV región", "Segunda", 100, 4000, 4500, 4250, "`$/caja 10 kilos", "Región de Arica y Parinacota", 425, 10, "Hortaliza")
  ,@(1, "Agrícola del Norte S.A. de Arica", "Arica y Parinacota", 44526, 15, 100112043, "Pepino dulce", "Cultivar XV región", "Tercera", 120, 3000, 3500, 3250, "`$/caja 10 kilos", "Región de Arica y Parinacota", 325, 10, "Hortaliza")
  ,@(1, "Agrícola del Norte S.A. de Arica", "Arica y Parinacota", 45021, 15, 100112043, "Pepino dulce", "Cultivar IV Región", "Primera", 270, 17000, 18000, 17500, "`$/bandeja 18 kilos", "Provincia de Limarí", 972, 18, "Hortaliza")
  ,@(1, "Agrícola del Norte S.A. de Arica", "Arica y Parinacota", 44221, 15, 100112043, "Pepino dulce", "Cultivar XV región", "Primera", 140, 5000, 6000, 5500, "`$/caja 10 kilos", "Región de Arica y Parinacota", 550, 10, "Hortaliza")
  ,@(1, "Agrícola del Norte S.A. de Arica", "Arica y Parinacota", 44377, 15, 100112043, "Pepino dulce", "Cultivar IV Región", "Primera", 100, 17000, 18000, 17600, "`$/bandeja 18 kilos", "Provincia de Limarí", 978, 18, "Hortaliza")
  ,@(1, "Agrícola del Norte S.A. de Arica", "Arica y Parinacota", 44533, 15, 100112043, "Pepino dulce", "Cultivar XV región", "Primera", 100, 6000, 7000, 6500, "`$/caja 10 kilos", "Región de Arica y Parinacota", 650, 10, "Hortaliza")
  ,@(1, "Agrícola del Norte S.A. de Arica", "Arica y Parinacota", 44533, 15, 100112043, "Pepino dulce", "Cultivar XV región", "Segunda", 120, 4000, 5000, 4500, "`$/caja 10 kilos", "Región de Arica y Parinacota", 450, 10, "Hortaliza")
  ,@(1, "Agrícola del Norte S.A. de Arica", "Arica y Parinacota", 45042, 15, 100112043, "Pepino dulce", "Cultivar IV Región", "Segunda", 220, 17000, 18000, 17545, "`$/bandeja 18 kilos", "Provincia de Limarí", 975, 18, "Hortaliza")
  ,@(1, "Agrícola del Norte S.A. de Arica", "Arica y Parinacota", 44391, 15, 100112043, "Pepino dulce", "Cultivar IV Región", "Segunda", 100, 15000, 16000, 15500, "`$/bandeja 18 kilos", "Provincia de Limarí", 861, 18, "Hortaliza")
  ,@(1, "Agrícola del Norte S.A. de Arica", "Arica y Parinacota", 44769, 15, 100112043, "Pepino dulce", "Cultivar IV Región", "Primera", 140, 17000, 18000, 17500, "`$/bandeja 18 kilos", "Provincia de Limarí", 972, 18, "Hortaliza")
  ,@(1, "Agrícola del Norte S.A. de Arica", "Arica y Parinacota", 44412, 15, 100112043, "Pepino dulce", "Cultivar IV Región", "Primera", 150, 17000, 18000, 17500, "`$/bandeja 18 kilos", "Provincia de Limarí", 972, 18, "Hortaliza")
  ,@(1, "Agrícola del Norte S.A. de Arica", "Arica y Parinacota", 44363, 15, 100112043, "Pepino dulce", "Cultivar IV Región", "Primera", 140, 14000, 15000, 14500, "`$/bandeja 18 kilos", "Provincia de Limarí", 806, 18, "Hortaliza")
  ,@(1, "Agrícola del Norte S.A. de Arica", "Arica y Parinacota", 45035, 15, 100112043, "Pepino dulce", "Cultivar IV Región", "Primera", 250, 19000, 20000, 19500, "`$/bandeja 18 kilos", "Provincia de Limarí", 1083, 18, "Hortaliza")
  ,@(1, "Agrícola del Norte S.A. de Arica", "Arica y Parinacota", 44776, 15, 100112043, "Pepino dulce", "Cultivar IV Región", "Primera", 200, 17000, 18000, 17500, "`$/bandeja 18 kilos", "Provincia de Limarí", 972, 18, "Hortaliza")
  ,@(1, "Agrícola del Norte S.A. de Arica", "Arica y Parinacota", 44554, 15, 100112043, "Pepino dulce", "Cultivar XV región", "Primera", 200, 5000, 6000, 5500, "`$/caja 10 kilos", "Región de Arica y Parinacota", 550, 10, "Hortaliza")
  ,@(1, "Agrícola del Norte S.A. de Arica", "Arica y Parinacota", 44783, 15, 100112043, "Pepino dulce", "Cultivar IV Región", "Primera", 150, 17000, 18000, 17500, "`$/bandeja 18 kilos", "Provincia de Limarí", 972, 18, "Hortaliza")
  ,@(1, "Agrícola del Norte S.A. de Arica", "Arica y Parinacota", 44398, 15, 100112043, "Pepino dulce", "Cultivar IV Región", "Primera", 100, 17000, 18000, 17500, "`$/bandeja 18 kilos", "Provincia de Limarí", 972, 18, "Hortaliza")
  ,@(1, "Agrícola del Norte S.A. de Arica", "Arica y Parinacota", 44398, 15, 100112043, "Pepino dulce", "Cultivar IV Región", "Segunda", 100, 15000, 16000, 15500, "`$/bandeja 18 kilos", "Provincia de Limarí", 861, 18, "Hortaliza")
  ,@(1, "Agrícola del Norte S.A. de Arica", "Arica y Parinacota", 44757, 15, 100112043, "Pepino dulce", "Cultivar XV región", "Primera", 150, 6000, 6500, 6250, "`$/caja 10 kilos", "Región de Arica y Parinacota", 625, 10, "Hortaliza")
  ,@(1, "Agrícola del Norte S.A. de Arica", "Arica y Parinacota", 44211, 15, 100112043, "Pepino dulce", "Cultivar XV región", "Segunda", 140, 4500, 5000, 4750, "`$/caja 10 kilos", "Región de Arica y Parinacota", 475, 10, "Hortaliza")
  ,@(1, "Agrícola del Norte S.A. de Arica", "Arica y Parinacota", 44762, 15, 100112043, "Pepino dulce", "Cultivar IV Región", "Primera", 160, 15000, 16000, 15500, "`$/bandeja 18 kilos", "Provincia de Limarí", 861, 18, "Hortaliza")
  ,@(1, "Agrícola del Norte S.A. de Arica", "Arica y Parinacota", 44454, 15, 100112043, "Pepino dulce", "Cultivar IV Región", "Primera", 160, 19000, 20000, 19500, "`$/bandeja 18 kilos", "Provincia de Limarí", 1083, 18, "Hortaliza")
)

for ($i = 0; $i -lt $data.Count; $i++) {
  $rowIndex = $i + 2
  $rowValues = $data[$i]
  for ($j = 0; $j -lt $rowValues.Count; $j++) {
    $colIndex = $j + 1
    $ws.Cells.Item($rowIndex, $colIndex).Value = $rowValues[$j]
  }
}

$ws.Range("D2:D33").NumberFormat = "YYYY-MM-DD HH:MM:SS"
